$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row for "30 days" payment reminder entry (Natuurlik Free Delivery option)
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "30 days"
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = $false

# Expand the table range to include the new row
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:D11"))

# Update selection to match target state
$ws.Range("F8").Select()
